# Update "想去人数" (people interested) counts in column F across sheets.
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 247
$ws.Range("F3").Value = 586
$ws.Range("F6").Value = 3189
$ws.Range("F7").Value = 2762
$ws.Range("F9").Value = 49
$ws.Range("F10").Value = 21
$ws.Range("F11").Value = 357
$ws.Range("F12").Value = 294
$ws.Range("F14").Value = 5753
$ws.Range("F15").Value = 619
$ws.Range("F17").Value = 62
$ws.Range("F20").Value = 470
$ws.Range("F21").Value = 1251
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = 1054
$ws.Range("F25").Value = 133
$ws.Range("F26").Value = 335
$ws.Range("F27").Value = 46

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F9").Value = 63
$ws.Range("F18").Value = 53
$ws.Range("F30").Value = 63
$ws.Range("F33").Value = 2

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F6").Value = 1120
$ws.Range("F9").Value = 1447
$ws.Range("F13").Value = 526

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F6").Value = 1120
$ws.Range("F7").Value = 1447
$ws.Range("F11").Value = 247
$ws.Range("F12").Value = 586
$ws.Range("F14").Value = 3189
$ws.Range("F15").Value = 2762
$ws.Range("F17").Value = 49
$ws.Range("F18").Value = 21
$ws.Range("F19").Value = 357
$ws.Range("F21").Value = 63
$ws.Range("F22").Value = 295
$ws.Range("F25").Value = 619
$ws.Range("F28").Value = 62
$ws.Range("F31").Value = 470
$ws.Range("F34").Value = 53
$ws.Range("F38").Value = 293
$ws.Range("F39").Value = 1251
$ws.Range("F42").Value = 63
$ws.Range("F43").Value = 1056
$ws.Range("F46").Value = 133
$ws.Range("F47").Value = 335
$ws.Range("F48").Value = 46
